$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.578.39'
$ws.Range("E2").Value = '  +0.15%  '

$ws.Range("D3").Value = '1.755.89'
$ws.Range("E3").Value = '  +0.24%  '

$ws.Range("E4").Value = '  -0.23%  '

$ws.Range("D5").Value = '324.26'
$ws.Range("E5").Value = '  +0.00%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.13%  '

$ws.Range("D7").Value = '0.4574'
$ws.Range("E7").Value = '  +2.05%  '

$ws.Range("D8").Value = '0.3570'
$ws.Range("E8").Value = '  -1.31%  '

$ws.Range("D9").Value = '0.07479'
$ws.Range("E9").Value = '  -0.33%  '

$ws.Range("D10").Value = '41.44'
$ws.Range("E10").Value = '  -1.70%  '

$ws.Range("D11").Value = '1.087'
$ws.Range("E11").Value = '  -1.59%  '

$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.10%  '

$ws.Range("D13").Value = '20.81'
$ws.Range("E13").Value = '  +0.64%  '

$ws.Range("D14").Value = '6.023'
$ws.Range("E14").Value = '  -0.17%  '

$ws.Range("D15").Value = '7.190'
$ws.Range("E15").Value = '  +0.11%  '

$ws.Range("D16").Value = '1.749.79'
$ws.Range("E16").Value = '  -0.11%  '

$ws.Range("D17").Value = '94.50'
$ws.Range("E17").Value = '  +1.84%  '

$ws.Range("D18").Value = '0.00001058'
$ws.Range("E18").Value = '  -0.64%  '

$ws.Range("D19").Value = '0.06396'
$ws.Range("E19").Value = '  -0.43%  '

$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.10%  '

$ws.Range("D21").Value = '17.15'
$ws.Range("E21").Value = '  +0.66%  '

$ws.Range("D22").Value = '5.759'
$ws.Range("E22").Value = '  -1.51%  '

$ws.Range("D23").Value = '27.620.91'

$ws.Range("E24").Value = '  -0.40%  '

$ws.Range("E25").Value = '  -0.90%  '

$ws.Range("D26").Value = '165.69'
$ws.Range("E26").Value = '  +1.83%  '

$ws.Range("D27").Value = '20.17'
$ws.Range("E27").Value = '  -1.34%  '

$ws.Range("D28").Value = '1.954.48'
$ws.Range("E28").Value = '  +0.17%  '

$ws.Range("D29").Value = '2.121'
$ws.Range("E29").Value = '  -0.88%  '

$ws.Range("D30").Value = '125.71'
$ws.Range("E30").Value = '  +0.16%  '

$ws.Range("D31").Value = '1.084'
$ws.Range("E31").Value = '  -0.20%  '

$ws.Range("D32").Value = '0.09214'
$ws.Range("E32").Value = '  +2.42%  '

$ws.Range("D33").Value = '3.653'
$ws.Range("E33").Value = '  +0.47%  '

$ws.Range("D34").Value = '5.530'
$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '11.75'
$ws.Range("E35").Value = '  -2.95%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.02287'
$ws.Range("E36").Value = '  -0.56%  '

$ws.Range("E37").Value = '  -0.32%  '

$ws.Range("D38").Value = '0.06021'
$ws.Range("E38").Value = '  +0.88%  '

$ws.Range("D39").Value = '0.6286'
$ws.Range("E39").Value = '  -1.04%  '

$ws.Range("D40").Value = '4.935'
$ws.Range("E40").Value = '  +0.02%  '

$ws.Range("E41").Value = '  -0.87%  '

$ws.Range("D42").Value = '1.389'
$ws.Range("E42").Value = '  -0.22%  '

$ws.Range("D43").Value = '7.792'
$ws.Range("E43").Value = '  +0.06%  '

$ws.Range("D44").Value = '13.18'
$ws.Range("E44").Value = '  -0.72%  '

$ws.Range("D45").Value = '3.715'
$ws.Range("E45").Value = '  +0.05%  '

$ws.Range("D46").Value = '0.5865'
$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").Value = '122.10'
$ws.Range("E47").Value = '  +0.42%  '

$ws.Range("D48").Value = '1.939'
$ws.Range("E48").Value = '  -0.79%  '

$ws.Range("D49").Value = '0.06895'
$ws.Range("E49").Value = '  +0.61%  '

$ws.Range("D50").Value = '1.131'
$ws.Range("E50").Value = '  -1.89%  '

$ws.Range("D51").Value = '71.77'
$ws.Range("E51").Value = '  -1.21%  '
